{"js": "// Apply the \"hearing-ready\" revision to 07_YTTRANDE_NY_BEVISNING.docx:\n//  - bump the letter date and the \"sista dag\" line to reflect the\n//    completed hearing (18 May 2026)\n//  - re-word the background bullets to past tense / completed status\n//  - drop the \"identity not yet revealed\" bullet, add the two new\n//    post-hearing bullets\n//  - rename the evidence heading (\"NY\" -> current) and update the\n//    translation + signature lines\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text,style\");\nawait context.sync();\n\n// ---- 1) simple text replacements (paragraph text is unique enough to\n//         key off directly) -------------------------------------------------\nconst textReplacements = [\n  [\"Datum: 2026-02-16\", \"Datum: 2026-05-18\"],\n  [\"Sista dag f\u00f6r ny bevisning: 2026-03-20\", \"Bevisning inl\u00e4mnad: 2026-03-20\"],\n  [\n    \"Vid muntlig f\u00f6rberedelse den 16 februari 2026 redovisades att:\",\n    \"Vid muntlig f\u00f6rberedelse den 16 februari 2026 beslutades att:\",\n  ],\n  [\n    \"Parterna har frist till den 20 mars 2026 att inkomma med ny bevisning\",\n    \"Parterna hade frist till den 20 mars 2026 att inkomma med ny bevisning\",\n  ],\n  [\"NY BEVISNING SOM \u00c5BEROPAS\", \"BEVISNING SOM \u00c5BEROPAS\"],\n  [\n    \"I enlighet med r\u00e4ttens anvisningar f\u00f6rses all arabisk bevisning med auktoriserad svensk \u00f6vers\u00e4ttning f\u00f6re fristen 2026-03-20.\",\n    \"All arabisk bevisning har f\u00f6rsetts med auktoriserad svensk \u00f6vers\u00e4ttning och inl\u00e4mnats till r\u00e4tten.\",\n  ],\n  [\"Eskilstuna den 2026-02-16\", \"Eskilstuna den 2026-05-18\"],\n];\n\n// paragraph to delete outright (identity of witnesses no longer \"not yet\n// revealed\" once the hearing has happened)\nconst textToDelete = \"Heba Alhussien har \u00e5beropat 2 vittnen vars identitet \u00e4nnu inte avsl\u00f6jats\";\n\n// anchor after which the two new post-hearing bullets get inserted\nconst insertAfterText = \"All arabisk bevisning ska f\u00f6rses med auktoriserad \u00f6vers\u00e4ttning\";\n\nlet deleteTarget = null;\nlet insertAfterTarget = null;\n\nfor (const p of paragraphs.items) {\n  const t = p.text;\n  for (const [oldText, newText] of textReplacements) {\n    if (t === oldText) {\n      p.insertText(newText, Word.InsertLocation.replace);\n    }\n  }\n  if (t === textToDelete) {\n    deleteTarget = p;\n  }\n  if (t === insertAfterText) {\n    insertAfterTarget = p;\n  }\n}\n\n// ---- 2) drop the stale \"identity not yet revealed\" bullet -----------------\nif (deleteTarget) {\n  deleteTarget.delete();\n}\n\n// ---- 3) add the two new bullets after \"All arabisk bevisning...\" ----------\nif (insertAfterTarget) {\n  const bullet2 = insertAfterTarget.insertParagraph(\n    \"Huvudf\u00f6rhandling best\u00e4mdes till den 18 maj 2026\",\n    Word.InsertLocation.after\n  );\n  bullet2.style = \"List Bullet\";\n\n  const bullet1 = insertAfterTarget.insertParagraph(\n    \"Heba Alhussien har \u00e5beropat 2 vittnen\",\n    Word.InsertLocation.after\n  );\n  bullet1.style = \"List Bullet\";\n}\n\nawait context.sync();\n", "ps1": "# Apply the \"hearing-ready\" revision to 07_YTTRANDE_NY_BEVISNING.docx:\n#  - bump the letter date and the \"sista dag\" line to reflect the\n#    completed hearing (18 May 2026)\n#  - re-word the background bullets to past tense / completed status\n#  - drop the \"identity not yet revealed\" bullet, add the two new\n#    post-hearing bullets\n#  - rename the evidence heading (\"NY\" -> current) and update the\n#    translation + signature lines\n\n$d = $word.ActiveDocument\n\nfunction Replace-DocText($findText, $replaceText) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\nfunction Find-ParagraphByText($text) {\n    $count = $d.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        $t = $p.Range.Text\n        if ($t.Length -gt 0) {\n            $t = $t.Substring(0, $t.Length - 1)\n        }\n        if ($t -eq $text) {\n            return $p\n        }\n    }\n    return $null\n}\n\n# ---- 1) simple text replacements ------------------------------------------\nReplace-DocText \"Datum: 2026-02-16\" \"Datum: 2026-05-18\"\nReplace-DocText \"Sista dag f\u00f6r ny bevisning: 2026-03-20\" \"Bevisning inl\u00e4mnad: 2026-03-20\"\nReplace-DocText \"Vid muntlig f\u00f6rberedelse den 16 februari 2026 redovisades att:\" \"Vid muntlig f\u00f6rberedelse den 16 februari 2026 beslutades att:\"\nReplace-DocText \"Parterna har frist till den 20 mars 2026 att inkomma med ny bevisning\" \"Parterna hade frist till den 20 mars 2026 att inkomma med ny bevisning\"\nReplace-DocText \"NY BEVISNING SOM \u00c5BEROPAS\" \"BEVISNING SOM \u00c5BEROPAS\"\nReplace-DocText \"I enlighet med r\u00e4ttens anvisningar f\u00f6rses all arabisk bevisning med auktoriserad svensk \u00f6vers\u00e4ttning f\u00f6re fristen 2026-03-20.\" \"All arabisk bevisning har f\u00f6rsetts med auktoriserad svensk \u00f6vers\u00e4ttning och inl\u00e4mnats till r\u00e4tten.\"\nReplace-DocText \"Eskilstuna den 2026-02-16\" \"Eskilstuna den 2026-05-18\"\n\n# ---- 2) drop the stale \"identity not yet revealed\" bullet -----------------\n$stale = Find-ParagraphByText \"Heba Alhussien har \u00e5beropat 2 vittnen vars identitet \u00e4nnu inte avsl\u00f6jats\"\nif ($stale -ne $null) {\n    $stale.Range.Delete()\n}\n\n# ---- 3) add the two new bullets after \"All arabisk bevisning...\" ----------\n$anchor = Find-ParagraphByText \"All arabisk bevisning ska f\u00f6rses med auktoriserad \u00f6vers\u00e4ttning\"\nif ($anchor -ne $null) {\n    $anchor.Range.InsertParagraphAfter()\n    $bullet1 = $anchor.Next()\n    $bullet1.Range.Text = \"Heba Alhussien har \u00e5beropat 2 vittnen\"\n    $bullet1.Style = \"List Bullet\"\n\n    $bullet1.Range.InsertParagraphAfter()\n    $bullet2 = $bullet1.Next()\n    $bullet2.Range.Text = \"Huvudf\u00f6rhandling best\u00e4mdes till den 18 maj 2026\"\n    $bullet2.Style = \"List Bullet\"\n}\n"}
